# Restore the hour-registration numbers for "week 3.2" (rows 16-19) and a
# couple of leftover zero entries in "week 3.3" (rows 21-22) that had been
# cleared out. The weekly total row (row 29) recalculates automatically via
# its existing SUM formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 16 (Tuesday, week 3.2)
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 5
$ws.Range("I16").Value = 0

# Row 17 (Wednesday, week 3.2)
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 4.5
$ws.Range("F17").Value = 4
$ws.Range("G17").Value = 3
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 4.5

# Row 18 (Thursday, week 3.2)
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 4

# Row 19 (Friday, week 3.2)
$ws.Range("D19").Value = 1.5
$ws.Range("H19").Value = 0

# Row 21 (Sunday, week 3.3)
$ws.Range("H21").Value = 0

# Row 22 (Monday, week 3.3)
$ws.Range("H22").Value = 0
